# Updated detailed-accuracy results for ChloeTing Video 5 after re-running the
# keyword extraction/classification with a video selected from the new
# dropdown menu. The Predicted Category (column C) and Confidence Level
# (column D) values below reflect the refreshed run's output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = 99.99998807907104
    "D3" = 100
    "C4" = 2
    "D4" = 57.00954794883728
    "C5" = 2
    "D5" = 98.27585816383362
    "D8" = 99.99991655349731
    "C9" = 1
    "D9" = 99.99948740005493
    "D10" = 99.99998807907104
    "D11" = 99.99984502792358
    "D12" = 99.99967813491821
    "D15" = 99.99905824661255
    "D16" = 99.99994039535522
    "D19" = 100
    "D24" = 99.99946355819702
    "D25" = 99.90358352661133
    "D26" = 98.4968900680542
    "D27" = 99.99998807907104
    "D28" = 99.99998807907104
    "C29" = 2
    "D29" = 99.99982118606567
    "D30" = 99.99912977218628
    "D31" = 99.73229765892029
    "D33" = 99.39715266227722
    "D34" = 99.99998807907104
    "D35" = 100
    "D37" = 99.99995231628418
    "D38" = 67.62106418609619
    "D40" = 99.99984502792358
    "D41" = 63.95573616027832
    "D42" = 99.99816417694092
    "D43" = 99.99972581863403
    "D45" = 99.99977350234985
    "D46" = 100
    "D47" = 99.9996542930603
    "D48" = 99.99998807907104
    "D49" = 100
    "D50" = 99.99833106994629
    "D51" = 100
    "D52" = 100
    "D53" = 99.97274279594421
    "D54" = 99.99996423721313
    "D55" = 99.97361302375793
    "D57" = 98.32385778427124
    "D58" = 99.99998807907104
    "D60" = 99.6665894985199
    "D61" = 100
    "D62" = 99.87700581550598
    "D63" = 100
    "D64" = 99.99998807907104
    "D65" = 99.99943971633911
    "D66" = 99.9966025352478
    "C68" = 1
    "D68" = 91.51563048362732
    "D70" = 99.99386072158813
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
